$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 corresponds to the "codevita_exam_solutions" group entry.
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = "2026-02-22T20:33:23.863554+00:00"
$ws.Range("H3").Value = 46
$ws.Range("L3").Value = "[486988, 486982, 487002, 487051, 487036, 487065, 487059, 487110, 487103, 487102, 487114, 487074, 487066, 487055, 487138, 487137, 487149, 487164, 487144, 487159, 487228, 487218, 487227, 487257, 487230, 487240, 487291, 487304, 487314, 487315, 487370, 487384, 487372, 487416, 487434, 487439, 487419, 487433, 487602, 487609, 487607, 487608, 487612, 487605, 487601, 487610]"
